# Update the two rows whose records changed (the former duplicate rows are
# replaced with new, distinct person records).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: was a duplicate of the "Jessie Marlowe" record -> now "Michelle Norton"
$ws.Range("A1").Value = "Michelle"
$ws.Range("B1").Value = "Norton"
$ws.Range("C1").Value = "Aperture Inc."
$ws.Range("D1").Value = "Scientist"
$ws.Range("E1").Value = "13 White Rabbit Street"
$ws.Range("F1").Value = "mnorton@aperture.us"
$ws.Range("G1").Value = 40731254562

# Row 7: was a duplicate of the "Doug Derrick" record -> now "John Smith"
$ws.Range("A7").Value = "John"
$ws.Range("B7").Value = "Smith"
$ws.Range("C7").Value = "IT Solutions"
$ws.Range("D7").Value = "Analyst"
$ws.Range("E7").Value = "98 North Road"
$ws.Range("F7").Value = "jsmith@itsolutions.co.uk"
$ws.Range("G7").Value = 40716543298
